$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "(Done Pss 1-3" -> "(Done Pss 1-4" in the TODO list, and move the
#    "_GoBack" bookmark from its old spot (after "...removed from some
#    pages") to the position right after the new "(Done Pss 1-4" text,
#    i.e. right before the closing ")" run -- mirroring how Word parks
#    _GoBack at the location of the most recent edit.
# ------------------------------------------------------------------

# Remove the existing _GoBack bookmark (Word keeps only one at a time).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the run to edit.
$target = $d.Content
$null = $target.Find.Execute("(Done Pss 1-3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Drop a throwaway guard bookmark immediately before the match so the
# host engine's run re-coalescing doesn't fold the previous ("...
# consistent with psalter. ") run into the one we are about to edit.
$guardPoint = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("zzGuard", $guardPoint)

# Update the TODO text.
$editRange = $d.Content
$null = $editRange.Find.Execute("(Done Pss 1-3", $true, $false, $false, $false, $false, $true, 1, $false, "(Done Pss 1-4", 2)

if ($d.Bookmarks.Exists("zzGuard")) {
    $d.Bookmarks("zzGuard").Delete()
}

# Find the freshly-edited text again so we can drop a collapsed
# _GoBack bookmark right after it (before the trailing ")").
$goBackRange = $d.Content
$null = $goBackRange.Find.Execute("(Done Pss 1-4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ------------------------------------------------------------------
# 2) Footer page-number field result: stale cached "I" -> "II" (the
#    roman-numeral footer on the second/continuation page of that
#    section, now that the section runs a bit longer).
# ------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers(1)
    if ($ftr.Exists) {
        $fr = $ftr.Range
        $null = $fr.Find.Execute("I", $true, $false, $false, $false, $false, $true, 1, $false, "II", 2)
    }
}
